$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 5).Value = 0.27

# Row 4
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 0.06

# Row 5
$ws.Cells.Item(5, 3).Value = 5
$ws.Cells.Item(5, 4).Value = 2
$ws.Cells.Item(5, 5).Value = 0.08

# Row 6
$ws.Cells.Item(6, 3).Value = 2
$ws.Cells.Item(6, 4).Value = 8
$ws.Cells.Item(6, 5).Value = 0

# Row 7
$ws.Cells.Item(7, 4).Value = 2
$ws.Cells.Item(7, 5).Value = 0.27

# Row 9
$ws.Cells.Item(9, 3).Value = 41
$ws.Cells.Item(9, 4).Value = 37
$ws.Cells.Item(9, 5).Value = 0.05

# Row 11
$ws.Cells.Item(11, 3).Value = 0
$ws.Cells.Item(11, 4).Value = 4
$ws.Cells.Item(11, 5).Value = 0

# Row 12
$ws.Cells.Item(12, 3).Value = 8
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 0

# Row 17
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 15

# Row 18
$ws.Cells.Item(18, 3).Value = 1
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 0.37

# Row 19
$ws.Cells.Item(19, 3).Value = 8
$ws.Cells.Item(19, 4).Value = 11
$ws.Cells.Item(19, 5).Value = 0.07000000000000001

# Row 22
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 5).Value = 0.27

# Row 25
$ws.Cells.Item(25, 3).Value = 4
$ws.Cells.Item(25, 4).Value = 4
$ws.Cells.Item(25, 5).Value = 0.2

# Row 26
$ws.Cells.Item(26, 4).Value = 0

# Row 29
$ws.Cells.Item(29, 3).Value = 1
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 0.37

# Row 30
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 1

# Row 31
$ws.Cells.Item(31, 3).Value = 2
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0.14

# Row 33
$ws.Cells.Item(33, 3).Value = 6
$ws.Cells.Item(33, 4).Value = 3
$ws.Cells.Item(33, 5).Value = 0.09

# Row 34
$ws.Cells.Item(34, 3).Value = 7
$ws.Cells.Item(34, 4).Value = 7
$ws.Cells.Item(34, 5).Value = 0.15

# Row 35
$ws.Cells.Item(35, 3).Value = 9
$ws.Cells.Item(35, 4).Value = 6
$ws.Cells.Item(35, 5).Value = 0.09

$wb.Save()
